$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (59 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 184.42857
$ws.Range("I9").Value = 21.125
$ws.Range("J9").Value = 402.16666
$ws.Range("K9").Value = 21.125
$ws.Range("L9").Value = 402.16666
$ws.Range("M9").Value = 147.875
$ws.Range("N9").Value = -740.16666
$ws.Range("H12").Value = 443.625
$ws.Range("I12").Value = 457.33334
$ws.Range("J12").Value = 402.5
$ws.Range("K12").Value = 457.33334
$ws.Range("L12").Value = 402.5
$ws.Range("M12").Value = -287.33334
$ws.Range("N12").Value = -742.5
$ws.Range("H17").Value = 2210.1177
$ws.Range("J17").Value = 2436.5715
$ws.Range("L17").Value = 7309.7145
$ws.Range("N17").Value = -7645.7145
$ws.Range("H53").Value = 512
$ws.Range("I53").Value = 187.5
$ws.Range("J53").Value = 566.0833
$ws.Range("K53").Value = 187.5
$ws.Range("L53").Value = 566.0833
$ws.Range("M53").Value = 449.5
$ws.Range("N53").Value = -1840.0833
$ws.Range("H75").Value = 300037000
$ws.Range("J75").Value = 300037000
$ws.Range("L75").Value = 300037000
$ws.Range("N75").Value = -300038872
$ws.Range("H78").Value = 300037000
$ws.Range("J78").Value = 300037000
$ws.Range("L78").Value = 900111000
$ws.Range("N78").Value = -900120360
$ws.Range("H115").Value = 4601.4287
$ws.Range("I115").Value = 3302.5
$ws.Range("K115").Value = 9907.5
$ws.Range("M115").Value = -8340.5
$ws.Range("H127").Value = 2578.4
$ws.Range("I127").Value = 1930.6666
$ws.Range("K127").Value = 5791.9998
$ws.Range("M127").Value = -831.9997999999996
$ws.Range("H132").Value = 4169695.2
$ws.Range("I132").Value = 3032.8096
$ws.Range("K132").Value = 9098.4288
$ws.Range("M132").Value = -6568.4288
$ws.Range("H135").Value = 1363.3334
$ws.Range("I135").Value = 1058.4
$ws.Range("J135").Value = 2888
$ws.Range("K135").Value = 9525.6
$ws.Range("L135").Value = 25992
$ws.Range("M135").Value = -6990.6
$ws.Range("N135").Value = -31062
$ws.Range("H138").Value = 1887222.8
$ws.Range("I138").Value = 1851.9231
$ws.Range("J138").Value = 2420045
$ws.Range("K138").Value = 5555.7693
$ws.Range("L138").Value = 7260135
$ws.Range("M138").Value = -415.7692999999999
$ws.Range("N138").Value = -7270415

# --- Sheet: ARM (20 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 74033.625
$ws.Range("I32").Value = 74033.625
$ws.Range("K32").Value = 74033.625
$ws.Range("M32").Value = -73746.625
$ws.Range("H45").Value = 2841.6365
$ws.Range("I45").Value = 2778.9443
$ws.Range("K45").Value = 2778.9443
$ws.Range("M45").Value = -2401.9443
$ws.Range("H61").Value = 627620.1
$ws.Range("I61").Value = 2593.6
$ws.Range("K61").Value = 2593.6
$ws.Range("M61").Value = -2381.6
$ws.Range("H110").Value = 1013.5714
$ws.Range("I110").Value = 1013.5714
$ws.Range("K110").Value = 1013.5714
$ws.Range("M110").Value = 1031.4286
$ws.Range("H136").Value = 627620.1
$ws.Range("I136").Value = 2593.6
$ws.Range("K136").Value = 7780.799999999999
$ws.Range("M136").Value = -5230.799999999999

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4862.5
$ws.Range("I94").Value = 3983.3333
$ws.Range("K94").Value = 3983.3333
$ws.Range("M94").Value = -3532.3333
$ws.Range("H132").Value = 199946.28
$ws.Range("J132").Value = 199946.28
$ws.Range("L132").Value = 199946.28
$ws.Range("N132").Value = -210066.28

# --- Sheet: CRP (29 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2736.1292
$ws.Range("I31").Value = 1221.5714
$ws.Range("J31").Value = 5916.7
$ws.Range("K31").Value = 1221.5714
$ws.Range("L31").Value = 5916.7
$ws.Range("M31").Value = -926.5714
$ws.Range("N31").Value = -6506.7
$ws.Range("H34").Value = 2736.1292
$ws.Range("I34").Value = 1221.5714
$ws.Range("J34").Value = 5916.7
$ws.Range("K34").Value = 1221.5714
$ws.Range("L34").Value = 5916.7
$ws.Range("M34").Value = -1019.5714
$ws.Range("N34").Value = -6320.7
$ws.Range("H52").Value = 87749.5
$ws.Range("J52").Value = 87749.5
$ws.Range("L52").Value = 87749.5
$ws.Range("N52").Value = -88337.5
$ws.Range("H96").Value = 57812
$ws.Range("J96").Value = 57812
$ws.Range("L96").Value = 57812
$ws.Range("N96").Value = -63304
$ws.Range("H132").Value = 2575.348
$ws.Range("I132").Value = 1811.0625
$ws.Range("J132").Value = 4322.2856
$ws.Range("K132").Value = 5433.1875
$ws.Range("L132").Value = 12966.8568
$ws.Range("M132").Value = -2903.1875
$ws.Range("N132").Value = -18026.8568

# --- Sheet: CUL (14 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 569.6
$ws.Range("I121").Value = 283
$ws.Range("J121").Value = 999.5
$ws.Range("K121").Value = 849
$ws.Range("L121").Value = 2998.5
$ws.Range("M121").Value = 461
$ws.Range("N121").Value = -5618.5
$ws.Range("H132").Value = 1631.6774
$ws.Range("I132").Value = 1163.7727
$ws.Range("J132").Value = 2775.4443
$ws.Range("K132").Value = 10473.9543
$ws.Range("L132").Value = 24978.9987
$ws.Range("M132").Value = -7943.954299999999
$ws.Range("N132").Value = -30038.9987

# --- Sheet: GSM (29 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 445.66666
$ws.Range("J9").Value = 785
$ws.Range("L9").Value = 785
$ws.Range("N9").Value = -1125
$ws.Range("H20").Value = 11100.8
$ws.Range("I20").Value = 6668.3335
$ws.Range("J20").Value = 17749.5
$ws.Range("K20").Value = 6668.3335
$ws.Range("L20").Value = 17749.5
$ws.Range("M20").Value = -6423.3335
$ws.Range("N20").Value = -18239.5
$ws.Range("H122").Value = 15410.667
$ws.Range("I122").Value = 19616.5
$ws.Range("K122").Value = 58849.5
$ws.Range("M122").Value = -56399.5
$ws.Range("H126").Value = 2164.3
$ws.Range("I126").Value = 1684.3077
$ws.Range("J126").Value = 3055.7144
$ws.Range("K126").Value = 5052.9231
$ws.Range("L126").Value = 9167.143199999999
$ws.Range("M126").Value = -2582.9231
$ws.Range("N126").Value = -14107.1432
$ws.Range("H132").Value = 24392006
$ws.Range("I132").Value = 31251400
$ws.Range("J132").Value = 3048.111
$ws.Range("K132").Value = 93754200
$ws.Range("L132").Value = 9144.332999999999
$ws.Range("M132").Value = -93751670
$ws.Range("N132").Value = -14204.333

# --- Sheet: LTW (41 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1452.15
$ws.Range("I16").Value = 1413.5555
$ws.Range("K16").Value = 1413.5555
$ws.Range("M16").Value = -1243.5555
$ws.Range("H22").Value = 2323.24
$ws.Range("I22").Value = 1439.8
$ws.Range("J22").Value = 2544.1
$ws.Range("K22").Value = 1439.8
$ws.Range("L22").Value = 2544.1
$ws.Range("M22").Value = -1144.8
$ws.Range("N22").Value = -3134.1
$ws.Range("H27").Value = 2323.24
$ws.Range("I27").Value = 1439.8
$ws.Range("J27").Value = 2544.1
$ws.Range("K27").Value = 1439.8
$ws.Range("L27").Value = 2544.1
$ws.Range("M27").Value = -1332.8
$ws.Range("N27").Value = -2758.1
$ws.Range("H55").Value = 1961.75
$ws.Range("I55").Value = 1318.125
$ws.Range("K55").Value = 1318.125
$ws.Range("M55").Value = -1145.125
$ws.Range("H93").Value = 3775.4783
$ws.Range("I93").Value = 1988.6086
$ws.Range("J93").Value = 5562.3477
$ws.Range("K93").Value = 1988.6086
$ws.Range("L93").Value = 5562.3477
$ws.Range("M93").Value = -740.6086
$ws.Range("N93").Value = -8058.3477
$ws.Range("H100").Value = 166673330
$ws.Range("J100").Value = 100008000
$ws.Range("L100").Value = 100008000
$ws.Range("N100").Value = -100009082
$ws.Range("H133").Value = 93998.336
$ws.Range("J133").Value = 93998.336
$ws.Range("L133").Value = 93998.336
$ws.Range("N133").Value = -99058.336
$ws.Range("H136").Value = 11868.158
$ws.Range("I136").Value = 5182.727
$ws.Range("K136").Value = 15548.181
$ws.Range("M136").Value = -12998.181

# --- Sheet: WVR (22 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5209.125
$ws.Range("J96").Value = 5667.143
$ws.Range("L96").Value = 5667.143
$ws.Range("N96").Value = -8413.143
$ws.Range("H107").Value = 41667184
$ws.Range("I107").Value = 505.83334
$ws.Range("J107").Value = 83333864
$ws.Range("K107").Value = 1517.50002
$ws.Range("L107").Value = 250001592
$ws.Range("M107").Value = 402.4999800000001
$ws.Range("N107").Value = -250005432
$ws.Range("H126").Value = 11907558
$ws.Range("I126").Value = 2375.75
$ws.Range("J126").Value = 27781134
$ws.Range("K126").Value = 7127.25
$ws.Range("L126").Value = 83343402
$ws.Range("M126").Value = -4657.25
$ws.Range("N126").Value = -83348342
$ws.Range("H132").Value = 12823021
$ws.Range("I132").Value = 15875489
$ws.Range("K132").Value = 47626467
$ws.Range("M132").Value = -47623937
